$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Title shape ("Testing" + " " + "custom" + " " + "properties") ---
# The flattened text is already "Testing custom properties", so a plain
# re-assignment is a no-op for the engine's run-diffing. Force a structural
# rewrite via a differing placeholder value first, then write the final text
# so it collapses to a single run.
$titleRange = $s.Shapes.Item(1).TextFrame.TextRange
$titleRange.Text = "PLACEHOLDER"
$titleRange.Text = "Testing custom properties"

# --- Subtitle shape (two leading line breaks, then "A." + " " + "M.") ---
# Only touch the run text after the two <a:br/> breaks, so the breaks
# themselves stay untouched, and use the same placeholder trick to force
# those runs to consolidate into a single run.
$subtitleRange = $s.Shapes.Item(2).TextFrame.TextRange
$breaksLen = 2
$runStart = $breaksLen + 1

$runsRange = $subtitleRange.Characters($runStart, $subtitleRange.Length - $breaksLen)
$runsRange.Text = "PLACEHOLDER"

$runsRange = $subtitleRange.Characters($runStart, $subtitleRange.Length - $breaksLen)
$runsRange.Text = "A. M."
